$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 101, shifting existing rows 101-129 down to 102-130.
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new weekly price observation.
$ws.Cells.Item(101, 1).Value = 6
$ws.Cells.Item(101, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(101, 3).Value = "Metropolitana"
$ws.Cells.Item(101, 4).Value = 44551
$ws.Cells.Item(101, 5).Value = 13
$ws.Cells.Item(101, 6).Value = "Fruta"
$ws.Cells.Item(101, 7).Value = 100101
$ws.Cells.Item(101, 8).Value = "Berries"
$ws.Cells.Item(101, 9).Value = 100101004
$ws.Cells.Item(101, 10).Value = "Frambuesa"
$ws.Cells.Item(101, 11).Value = "Sin especificar"
$ws.Cells.Item(101, 12).Value = "Primera"
$ws.Cells.Item(101, 13).Value = 500
$ws.Cells.Item(101, 14).Value = 8000
$ws.Cells.Item(101, 15).Value = 8000
$ws.Cells.Item(101, 16).Value = 8000
$ws.Cells.Item(101, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(101, 18).Value = "Provincia de Linares"
$ws.Cells.Item(101, 19).Value = 4000
$ws.Cells.Item(101, 20).Value = 2
